$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New CRM test rows (rows 3-7 have full data, rows 8-10 only have C/D/E -
# the batch/CRM values were never entered for that 3rd test run).
$ws.Range("A3").Value = 43173
$ws.Range("A4").Value = 43173
$ws.Range("A5").Value = 43173
$ws.Range("A6").Value = 43173
$ws.Range("A7").Value = 43173

# Match the date formatting already used on A2 (built-in short-date style)
# instead of minting a new number format.
$ws.Range("A2").Copy()
$ws.Range("A3:A7").PasteSpecial(-4122)

$ws.Range("B3").Value = 2217.5163287970099
$ws.Range("B4").Value = 2219.6900936089201
$ws.Range("B5").Value = 2225.21466328611
$ws.Range("B6").Value = 2221.4246394258698
$ws.Range("B7").Value = 2227.6089961644998

$ws.Range("C3").Value = 2207.0300000000002
$ws.Range("C4").Value = 2207.0300000000002
$ws.Range("C5").Value = 2207.0300000000002
$ws.Range("C6").Value = 2207.0300000000002
$ws.Range("C7").Value = 2207.0300000000002
$ws.Range("C8").Value = 2207.0300000000002
$ws.Range("C9").Value = 2207.0300000000002
$ws.Range("C10").Value = 2207.0300000000002

$ws.Range("D3:D10").Formula = "=100*(B3-C3)/C3"

$ws.Range("E3").Value = 169
$ws.Range("E4").Value = 169
$ws.Range("E5").Value = 169
$ws.Range("E6").Value = 169
$ws.Range("E7").Value = 169
$ws.Range("E8").Value = 169
$ws.Range("E9").Value = 169
$ws.Range("E10").Value = 169

$ws.Range("A9").Select()
